# VerveStacks CHE grids model update - 2025-08-31 21:16
# - Removes the re_profiles and load_shapes sheets (data moved elsewhere / no longer needed here)
# - Refreshes the "daynight" timeslice-grouping lookup table on ev_charging_uc so that
#   the Day/Night buckets point at the full explicit timeslice lists instead of the
#   single catch-all "AllSaAllH" placeholder.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the two obsolete sheets -----------------------------------------
$wb.Worksheets("re_profiles").Delete()
$wb.Worksheets("load_shapes").Delete()

# --- Refresh the lookup table on ev_charging_uc -----------------------------
$ws1 = $wb.Worksheets("ev_charging_uc")

$dayTimeslices = "RaD,SaD,TaD,Tc0916h07,Tc0916h08,Tc0916h09,Tc0916h10,Tc0916h11,Tc0916h12,Tc0916h13,Tc0916h14,Tc0916h15,Tc0916h16,Tc0916h17,Tc0916h18,Tc0917h07,Tc0917h08,Tc0917h09,Tc0917h10,Tc0917h11,Tc0917h12,Tc0917h13,Tc0917h14,Tc0917h15,Tc0917h16,Tc0917h17,Tc0917h18,WaD,Wb0121h07,Wb0121h08,Wb0121h09,Wb0121h10,Wb0121h11,Wb0121h12,Wb0121h13,Wb0121h14,Wb0121h15,Wb0121h16,Wb0121h17,Wb0121h18,Wb0122h07,Wb0122h08,Wb0122h09,Wb0122h10,Wb0122h11,Wb0122h12,Wb0122h13,Wb0122h14,Wb0122h15,Wb0122h16,Wb0122h17,Wb0122h18"
$nightTimeslices = "RaD,RaN,RaP,SaD,SaN,SaP,TaD,TaN,TaP,Tc0916h01,Tc0916h02,Tc0916h03,Tc0916h04,Tc0916h05,Tc0916h06,Tc0916h19,Tc0916h20,Tc0916h21,Tc0916h22,Tc0916h23,Tc0916h24,Tc0917h01,Tc0917h02,Tc0917h03,Tc0917h04,Tc0917h05,Tc0917h06,Tc0917h19,Tc0917h20,Tc0917h21,Tc0917h22,Tc0917h23,Tc0917h24,WaD,WaN,WaP,Wb0121h01,Wb0121h02,Wb0121h03,Wb0121h04,Wb0121h05,Wb0121h06,Wb0121h19,Wb0121h20,Wb0121h21,Wb0121h22,Wb0121h23,Wb0121h24,Wb0122h01,Wb0122h02,Wb0122h03,Wb0122h04,Wb0122h05,Wb0122h06,Wb0122h19,Wb0122h20,Wb0122h21,Wb0122h22,Wb0122h23,Wb0122h24"

$ws1.Range("B12").Value = "daynight"
$ws1.Range("C12").Value = "TimeSlices"

$ws1.Range("B13").Value = "Day"
$ws1.Range("C13").Value = $dayTimeslices

$ws1.Range("B14").Value = "Night"
$ws1.Range("C14").Value = $nightTimeslices
